# Add new columns I (I0) and J (IF) with header + data values, mirroring
# the styling and structure of the existing column H (IP).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the existing header cell H1 (bold, centered, bordered)
# onto the new header cells so they keep the same look / reuse the style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-42 (column I = "I0", column J = "IF")
$iValues = @(8,8,8,4,6,3,6,8,8,3,7,6,5,7,5,6,8,9,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1)
$jValues = @(9,9,9,7,7,7,7,8,8,5,9,6,7,8,7,9,9,9,6,6,5,2,6,7,6,4,7,5,7,7,5,4,6,7,6,6,5,6,5,3,3)

for ($r = 2; $r -le 42; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
